# "Tracing fix and removed screenshot"
# Update the sample "items" search-results sheet: swap the IFB washing-machine
# row for a Samsung AC row, fix the recorded colour for the iPhone row, make
# the header row's fill consistent across all four columns, and leave the
# selection on the last edited cell (D3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

# Header row (Result2 / Result3) picked up the same bold + yellow-fill look
# that SearchKey/Result1 already have.
$ws.Range("C1:D1").Interior.Color = 65535

# Row 2 used to be a hyperlink-style search result ("ifb washing machine").
# Replace it with a fresh, plain (non-hyperlink) Samsung AC entry.
$ws.Range("A2:B2").Style = "Normal"
$ws.Range("A2").NumberFormat = "0"

$ws.Range("A2").Value = "samsung 1.5 ton 5star"
$ws.Range("B2").Value = "Samsung"
$ws.Range("C2").Value = "1.5 Ton"

# Row 3 (iPhone 13 mini) keeps its data/formatting - only the recorded
# colour result was wrong.
$ws.Range("D3").Value = "Pink"

# Leave the cursor on the cell that was last fixed.
$ws.Range("D3").Select()
